# Fruta / hortaliza, semanal
# Insert a new weekly record as row 67 (pushing the existing rows 67-106
# down to 68-107) in the "Vega Central Mapocho de Santiago - Frambuesa"
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 67..106 down by one to make room for the new record.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new weekly data point.
$ws.Range("A67").Value = 9
$ws.Range("B67").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = 44673
$ws.Range("E67").Value = 13
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100101
$ws.Range("H67").Value = "Berries"
$ws.Range("I67").Value = 100101004
$ws.Range("J67").Value = "Frambuesa"
$ws.Range("K67").Value = "Sin especificar"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 250
$ws.Range("N67").Value = 8000
$ws.Range("O67").Value = 8000
$ws.Range("P67").Value = 8000
$ws.Range("Q67").Value = "$/bandeja 2 kilos"
$ws.Range("R67").Value = "Provincia de Linares"
$ws.Range("S67").Value = 4000
$ws.Range("T67").Value = 2
